$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.01179113890272636, 0.9527065656392815, 4.106302077608952, 2.026401262733754, 2.046136859441977, 52)
    3  = @(0.06492620600918006, 1.070932841330739, 4.414872682735234, 2.101159842262181, 2.121054081328115, 51)
    4  = @(0.0188834649839998, 0.9530525523764629, 3.865729941648752, 1.966145961430319, 1.986015741865533, 50)
    5  = @(0.09692881375077184, 1.09440480910644, 4.629837775659457, 2.151705782782455, 2.171796906056384, 49)
    6  = @(-0.000755637588625957, 0.9503968899563854, 3.951250381626263, 1.987775234181738, 2.0088103333927, 48)
    7  = @(0.03928990789677742, 1.033046348395609, 4.93222755344469, 2.220861894275439, 2.252012552851068, 36)
    8  = @(0.03868193857282701, 1.016268278851836, 5.049011224985318, 2.247000495101262, 2.279467318072585, 35)
    9  = @(0.02578329123745685, 1.456619164208787, 9.170609166841203, 3.028301366581801, 3.115983453452396, 18)
    10 = @(-0.8477958937427421, 1.04746007141034, 6.235465529401324, 2.497091413905651, 2.463407886881444, 11)
    11 = @(-0.248566666372188, 0.301374533835348, 0.1561541005302556, 0.3951633846021865, 0.3434543508584382, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i  # column B is 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
